$d = $word.ActiveDocument

# Locate the run containing "19.12.2023. godine" (the combo_dani day placeholder)
$rng = $d.Content
$found = $rng.Find.Execute("19.12.2023. godine", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start

    # Sub-range covering just the "19" portion that becomes the {{ dan_z }} placeholder
    $dayRange = $d.Range($start, $start + 2)
    $dayRange.Text = "{{ dan_z }}"

    # Toggling formatting on the replaced sub-range forces it to stay a distinct run
    # from the untouched ".12.2023. godine" text that follows it, even though both
    # runs end up with identical formatting.
    $dayRange.Font.Bold = $true
    $dayRange.Font.Bold = $false
}
